# Battle sounds added to Audio List: card attack sounds, Lion Death, card selection sound.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 22: sliding rock tablet (card selection sound) ---
$ws.Range("B22").Value = "sliding rock tablet"
$ws.Range("C22").Value = "SFX"
$ws.Range("D22").Value = "open"
$ws.Range("E22").Value = "yes"
$ws.Range("G22").Value = "each attack appears as a stone tablet, when you select it (or mouse over?) it slides up with a ""sliding rocks"" noise."

# --- Row 23: Horned Lion death ---
$ws.Range("B23").Value = "Horned Lion death"
$ws.Range("C23").Value = "SFX"
$ws.Range("D23").Value = "open"
$ws.Range("E23").Value = "yes"

# --- Row 24: Healing Magic sound ---
$ws.Range("B24").Value = "Healing Magic sound"
$ws.Range("C24").Value = "SFX"
$ws.Range("D24").Value = "open"
$ws.Range("E24").Value = "yes"

# --- Row 25: Magic shield sound (new row, A033) ---
$ws.Range("B25").Value = "Magic shield sound"
$ws.Range("A25").Value = "A033"
$ws.Range("C25").Value = "SFX"
$ws.Range("D25").Value = "open"
$ws.Range("E25").Value = "yes"

# --- Row 26: Lound swinging sound (new row, A034) ---
$ws.Range("B26").Value = "Lound swinging sound"
$ws.Range("G26").Value = "for powerful attacks "
$ws.Range("A26").Value = "A034"
$ws.Range("C26").Value = "SFX"
$ws.Range("D26").Value = "open"
$ws.Range("E26").Value = "yes"

# Update selection to match author's last active cell
$ws.Range("C29").Select()
